$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "582.24", "1.00") are stored as literal text, matching the source data,
# not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "69.317.22"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "3.427.72"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "582.24"
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("D6").Value = "178.33"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  +6.25%  "
$ws.Range("E10").Value = "  +0.42%  "
$ws.Range("D11").Value = "48.17"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").Value = "679.13"
$ws.Range("E13").Value = "  -2.30%  "
$ws.Range("D14").Value = "3.972.87"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").Value = "8.65"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "69.445.01"
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").Value = "3.418.09"
$ws.Range("E17").Value = "  +1.14%  "
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").Value = "17.78"
$ws.Range("E19").Value = "  +1.43%  "
$ws.Range("D20").Value = "11.32"
$ws.Range("E20").Value = "  +0.49%  "
$ws.Range("D21").Value = "0.911"
$ws.Range("E21").Value = "  +1.46%  "
$ws.Range("E22").Value = "  -2.88%  "
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").Value = "100.72"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "2.70"
$ws.Range("E26").Value = "  -0.50%  "
$ws.Range("D27").Value = "9.66"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").Value = "33.63"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("D29").Value = "8.76"
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("D30").Value = "6.88"
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("D31").Value = "3.73"
$ws.Range("E31").Value = "  +10.05%  "
$ws.Range("D32").Value = "559.81"
$ws.Range("E32").Value = "  +1.76%  "
$ws.Range("D33").Value = "11.02"
$ws.Range("E33").Value = "  -0.98%  "
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "58.08"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("E36").Value = "  +0.23%  "
$ws.Range("D37").Value = "3.612.14"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("D39").Value = "34.98"
$ws.Range("E39").Value = "  +0.82%  "
$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  +8.81%  "
$ws.Range("E41").Value = "  +2.37%  "
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("D43").Value = "3.41"
$ws.Range("E43").Value = "  +5.51%  "
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("E46").Value = "  +0.46%  "
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("D48").Value = "1.41"
$ws.Range("E48").Value = "  +4.66%  "
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "130.95"
$ws.Range("E50").Value = "  -0.83%  "
$ws.Range("D51").Value = "2.70"
$ws.Range("E51").Value = "  +2.70%  "

# Revert the explicit Text style applied above so the cell style index
# matches the original workbook (no residual "s" attribute on D cells).
$ws.Range("D2:D51").Style = "Normal"
